# Add Airflow DAG pipeline_jours_bien_etre
# This script nudges the position/size of several pictures, text boxes and
# connector arrows on slide 1 so a newly-added DAG box ("pipeline_jours_bien_etre")
# fits into the existing Airflow diagram.
#
# NOTE on numeric literals below: PowerPoint's Shape.Left/Top/Width/Height are
# single-precision (float32) values measured in points (1 pt = 12700 EMU).
# The host truncates (floors) the EMU value computed from the float32
# representation of whatever we assign, so plain "EMU/12700" literals can
# truncate to one EMU less than intended. The literals used here were solved
# so that floor(float32(literal) * 12700) lands exactly on the target EMU
# value from the target OOXML.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

# Picture 8 (right-hand icon, id=11): slide right a bit to make room
$sh = Get-ShapeById $s 11
$sh.Left = 712.1983464566929

# Connecteur droit avec flèche 15 (id=16): move start left
$sh = Get-ShapeById $s 16
$sh.Left = 414.1540157480315

# ZoneTexte 24 (id=25): shift left
$sh = Get-ShapeById $s 25
$sh.Left = 536.728058496063

# Connecteur droit avec flèche 25 (id=26): extend width
$sh = Get-ShapeById $s 26
$sh.Width = 75.79700887401575

# ZoneTexte 26 (id=27): shift right
$sh = Get-ShapeById $s 27
$sh.Left = 703.4259952519685

# Connecteur droit avec flèche 1047 (id=1048): resize to reach the new box
$sh = Get-ShapeById $s 1048
$sh.Width = 133.33047244094487
$sh.Height = 74.50275890551181

# Connecteur droit avec flèche 1049 (id=1050): resize to reach the new box
$sh = Get-ShapeById $s 1050
$sh.Width = 132.81165354330707
$sh.Height = 72.31456692913386

# ZoneTexte 1062 (id=1063): reposition slightly
$sh = Get-ShapeById $s 1063
$sh.Left = 594.2816535433071
$sh.Top = 449.60615573228347

# Image 1069 (id=1070): reposition and enlarge
$sh = Get-ShapeById $s 1070
$sh.Left = 631.4774015748031
$sh.Width = 75.26716535433071
$sh.Height = 19.733700787401574

# Picture 8 (left-hand icon, id=1108): slide right a bit
$sh = Get-ShapeById $s 1108
$sh.Left = 355.2513385826772

# ZoneTexte 1109 (id=1110): reposition slightly
$sh = Get-ShapeById $s 1110
$sh.Left = 337.76142932283466
$sh.Top = 146.25913385826772
